$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.044109778800644
$ws.Range("D2").Value = 1.051547795521744
$ws.Range("E2").Value = 1.05165889873041
$ws.Range("F2").Value = 1.06161246174445
$ws.Range("I2").Value = 1.035039329122346
$ws.Range("J2").Value = 1.049176623706213
$ws.Range("K2").Value = 1.054298993529602
$ws.Range("L2").Value = 1.054409788818577
$ws.Range("M2").Value = 1.064336048967985
$ws.Range("N2").Value = 1.050666574854116
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.045088906102785
$ws.Range("D3").Value = 1.052427921328022
$ws.Range("E3").Value = 1.052536067114041
$ws.Range("F3").Value = 1.062587453226967
$ws.Range("I3").Value = 1.035146189136606
$ws.Range("J3").Value = 1.049802725543827
$ws.Range("K3").Value = 1.054991614008665
$ws.Range("L3").Value = 1.055099481746786
$ws.Range("M3").Value = 1.065125291556436
$ws.Range("N3").Value = 1.051293565828155
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.045722916515794
$ws.Range("D4").Value = 1.052998172243164
$ws.Range("E4").Value = 1.053104456201761
$ws.Range("F4").Value = 1.063219344395726
$ws.Range("I4").Value = 1.035213896041238
$ws.Range("J4").Value = 1.050207686451972
$ws.Range("K4").Value = 1.055439896920143
$ws.Range("L4").Value = 1.055545921379931
$ws.Range("M4").Value = 1.065636368521232
$ws.Range("N4").Value = 1.051699101827239
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.045989561196936
$ws.Range("D5").Value = 1.053238084134265
$ws.Range("E5").Value = 1.053343597781327
$ws.Range("F5").Value = 1.063485231072236
$ws.Range("I5").Value = 1.03524201508004
$ws.Range("J5").Value = 1.050377890610572
$ws.Range("K5").Value = 1.055628380600714
$ws.Range("L5").Value = 1.055733642368375
$ws.Range("M5").Value = 1.065851316106283
$ws.Range("N5").Value = 1.051869547695268
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.046034338242518
$ws.Range("D6").Value = 1.05327837682037
$ws.Range("E6").Value = 1.05338376185575
$ws.Range("F6").Value = 1.063529888628663
$ws.Range("I6").Value = 1.0352467161475
$ws.Range("J6").Value = 1.050406466180987
$ws.Range("K6").Value = 1.055660029307075
$ws.Range("L6").Value = 1.055765163732146
$ws.Range("M6").Value = 1.06588741204426
$ws.Range("N6").Value = 1.051898163846269
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.045726479016955
$ws.Range("D7").Value = 1.053001377258218
$ws.Range("E7").Value = 1.053107650875156
$ws.Range("F7").Value = 1.063222896246594
$ws.Range("I7").Value = 1.035214273125258
$ws.Range("J7").Value = 1.050209960892209
$ws.Range("K7").Value = 1.055442415349272
$ws.Range("L7").Value = 1.055548429569903
$ws.Range("M7").Value = 1.065639240305941
$ws.Range("N7").Value = 1.051701379497442
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.044440585811109
$ws.Range("D8").Value = 1.051845081943805
$ws.Range("E8").Value = 1.051955174943943
$ws.Range("F8").Value = 1.061941755484674
$ws.Range("I8").Value = 1.035075740535943
$ws.Range("J8").Value = 1.049388252122206
$ws.Range("K8").Value = 1.054533044209477
$ws.Range("L8").Value = 1.054642839551593
$ws.Range("M8").Value = 1.064602696941359
$ws.Range("N8").Value = 1.050878503806734
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.042178161752513
$ws.Range("D9").Value = 1.049813348635986
$ws.Range("E9").Value = 1.049930569035917
$ws.Range("F9").Value = 1.05969198872737
$ws.Range("I9").Value = 1.034820633005201
$ws.Range("J9").Value = 1.047939046552157
$ws.Range("K9").Value = 1.052931514279355
$ws.Range("L9").Value = 1.053048363267131
$ws.Range("M9").Value = 1.062779168954785
$ws.Range("N9").Value = 1.049427240198549
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.040672275131688
$ws.Range("D10").Value = 1.048462840650646
$ws.Range("E10").Value = 1.048585080056355
$ws.Range("F10").Value = 1.058197444165698
$ws.Range("I10").Value = 1.034643199459826
$ws.Range("J10").Value = 1.04697212700807
$ws.Range("K10").Value = 1.051864494046498
$ws.Range("L10").Value = 1.051986308121639
$ws.Range("M10").Value = 1.061565564520412
$ws.Range("N10").Value = 1.048458947517803
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.040020788735786
$ws.Range("D11").Value = 1.047879015902738
$ws.Range("E11").Value = 1.048003492109274
$ws.Range("F11").Value = 1.057551563060962
$ws.Range("I11").Value = 1.034564629031248
$ws.Range("J11").Value = 1.046553266483234
$ws.Range("K11").Value = 1.051402633474813
$ws.Range("L11").Value = 1.051526659770704
$ws.Range("M11").Value = 1.061040568016095
$ws.Range("N11").Value = 1.048039492162978
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.039778884495385
$ws.Range("D12").Value = 1.047662301995003
$ws.Range("E12").Value = 1.047787618604117
$ws.Range("F12").Value = 1.057311845612988
$ws.Range("I12").Value = 1.034535183465061
$ws.Range("J12").Value = 1.046397657217155
$ws.Range("K12").Value = 1.051231103874398
$ws.Range("L12").Value = 1.051355961248561
$ws.Range("M12").Value = 1.060845637490895
$ws.Range("N12").Value = 1.047883661913892
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.039830769834512
$ws.Range("D13").Value = 1.047708781300365
$ws.Range("E13").Value = 1.047833917207735
$ws.Range("F13").Value = 1.057363257130822
$ws.Range("I13").Value = 1.03454151144246
$ws.Range("J13").Value = 1.046431037086427
$ws.Range("K13").Value = 1.051267896364133
$ws.Range("L13").Value = 1.051392575045578
$ws.Range("M13").Value = 1.06088744726847
$ws.Range("N13").Value = 1.047917089186408
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.040000791080115
$ws.Range("D14").Value = 1.04786109930651
$ws.Range("E14").Value = 1.047985644786009
$ws.Range("F14").Value = 1.057531744026992
$ws.Range("I14").Value = 1.034562200373295
$ws.Range("J14").Value = 1.046540404294084
$ws.Range("K14").Value = 1.051388454248721
$ws.Range("L14").Value = 1.051512549050839
$ws.Range("M14").Value = 1.061024453430077
$ws.Range("N14").Value = 1.048026611708044
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.040105558309036
$ws.Range("D15").Value = 1.047954966655396
$ws.Range("E15").Value = 1.048079149619224
$ws.Range("F15").Value = 1.05763557978494
$ws.Range("I15").Value = 1.034574912934464
$ws.Range("J15").Value = 1.046607785646851
$ws.Range("K15").Value = 1.051462737413075
$ws.Range("L15").Value = 1.0515864737083
$ws.Range("M15").Value = 1.061108877624778
$ws.Range("N15").Value = 1.04809408875006
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.040715524268071
$ws.Range("D16").Value = 1.048501607385522
$ws.Range("E16").Value = 1.048623699681597
$ws.Range("F16").Value = 1.058240335960382
$ws.Range("I16").Value = 1.034648377290744
$ws.Range("J16").Value = 1.046999921730454
$ws.Range("K16").Value = 1.05189514979967
$ws.Range("L16").Value = 1.052016818369736
$ws.Range("M16").Value = 1.06160041748439
$ws.Range("N16").Value = 1.048486781711881
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.04109829363055
$ws.Range("D17").Value = 1.048844757002942
$ws.Range("E17").Value = 1.048965554895923
$ws.Range("F17").Value = 1.058620023517924
$ws.Range("I17").Value = 1.034693993892556
$ws.Range("J17").Value = 1.047245851335698
$ws.Range("K17").Value = 1.052166435880118
$ws.Range("L17").Value = 1.052286824101944
$ws.Range("M17").Value = 1.061908882686828
$ws.Range("N17").Value = 1.048733060565373
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.041321611702072
$ws.Range("D18").Value = 1.049045002344034
$ws.Range("E18").Value = 1.04916505136349
$ws.Range("F18").Value = 1.058841610965378
$ws.Range("I18").Value = 1.034720433315938
$ws.Range("J18").Value = 1.047389280679076
$ws.Range("K18").Value = 1.052324688432343
$ws.Range("L18").Value = 1.052444335860228
$ws.Range("M18").Value = 1.062088853604218
$ws.Range("N18").Value = 1.04887669359487
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.04139776672802
$ws.Range("D19").Value = 1.049113296392988
$ws.Range("E19").Value = 1.049233091105378
$ws.Range("F19").Value = 1.0589171872336
$ws.Range("I19").Value = 1.034729419970914
$ws.Range("J19").Value = 1.047438183427301
$ws.Range("K19").Value = 1.052378651141461
$ws.Range("L19").Value = 1.052498046997403
$ws.Range("M19").Value = 1.062150227209489
$ws.Range("N19").Value = 1.048925665790607
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.041057220387189
$ws.Range("D20").Value = 1.048807930774059
$ws.Range("E20").Value = 1.048928866910132
$ws.Range("F20").Value = 1.058579273995808
$ws.Range("I20").Value = 1.034689117033172
$ws.Range("J20").Value = 1.047219467192469
$ws.Range("K20").Value = 1.052137327791831
$ws.Range("L20").Value = 1.052257852761104
$ws.Range("M20").Value = 1.061875782277382
$ws.Range("N20").Value = 1.048706638953635
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.039950721657914
$ws.Range("D21").Value = 1.04781624144844
$ws.Range("E21").Value = 1.047940960530388
$ws.Range("F21").Value = 1.057482123537534
$ws.Range("I21").Value = 1.034556115200099
$ws.Range("J21").Value = 1.046508199082608
$ws.Range("K21").Value = 1.051352952231269
$ws.Range("L21").Value = 1.051477218712884
$ws.Range("M21").Value = 1.060984106408235
$ws.Range("N21").Value = 1.047994360761474
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.039255524676662
$ws.Range("D22").Value = 1.047193564231871
$ws.Range("E22").Value = 1.047320717089739
$ws.Range("F22").Value = 1.056793410388859
$ws.Range("I22").Value = 1.03447098178193
$ws.Range("J22").Value = 1.046060847677795
$ws.Range("K22").Value = 1.050859934997094
$ws.Range("L22").Value = 1.050986608075908
$ws.Range("M22").Value = 1.060423917849645
$ws.Range("N22").Value = 1.047546374066354
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.039624013829161
$ws.Range("D23").Value = 1.04752357742842
$ws.Range("E23").Value = 1.047649434851925
$ws.Range("F23").Value = 1.057158404726424
$ws.Range("I23").Value = 1.034516255575535
$ws.Range("J23").Value = 1.046298010928711
$ws.Range("K23").Value = 1.051121278129654
$ws.Range("L23").Value = 1.051246670291935
$ws.Range("M23").Value = 1.060720842002538
$ws.Range("N23").Value = 1.04778387411629
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.041075779447415
$ws.Range("D24").Value = 1.048824570677647
$ws.Range("E24").Value = 1.048945444328679
$ws.Range("F24").Value = 1.058597686573587
$ws.Range("I24").Value = 1.034691321195043
$ws.Range("J24").Value = 1.04723138910442
$ws.Range("K24").Value = 1.052150480434714
$ws.Range("L24").Value = 1.052270943595129
$ws.Range("M24").Value = 1.061890738778209
$ws.Range("N24").Value = 1.048718577796068
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.042762634215133
$ws.Range("D25").Value = 1.050337904358363
$ws.Range("E25").Value = 1.050453235049892
$ws.Range("F25").Value = 1.06027267891675
$ws.Range("I25").Value = 1.034887884253527
$ws.Range("J25").Value = 1.048313843081121
$ws.Range("K25").Value = 1.053345435277081
$ws.Range("L25").Value = 1.053460414424799
$ws.Range("M25").Value = 1.063250232859562
$ws.Range("N25").Value = 1.049802568981574
